$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 117's "Numéro de page" value no longer applies -> becomes an empty
# text cell (same shape as every other day with nothing to report).
# A lone leading apostrophe forces text-entry semantics so the result stays
# an empty STRING (not a cleared/blank cell); resetting the Style afterwards
# drops the transient quote-prefix formatting it picks up along the way.
$ws.Range("C117").Value = "'"
$ws.Range("C117").Style = $ws.Range("D117").Style

# Append the new day's row produced by the latest script run.
$ws.Range("A118").Value = "'2025-05-26"
$ws.Range("A118").Style = $ws.Range("A117").Style
$ws.Range("B118").Value = "Rien ne nous concerne aujourd'hui !"
$ws.Range("C118").Value = "NA"
$ws.Range("D118").Value = 1
